# Apply "Created admin login functionality" changes to the improvements tracker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: new "Create admin login functionality" entry ------------------
$ws.Range("A21").Value = "ALL"
$ws.Range("B21").Value = "Create admin login functionality"
$ws.Range("C21").Value = "Added to login_tools and login_action, creating an 'admin' key for the `$_SESSION global"
$ws.Range("D21").Value = $true
$ws.Rows.Item(21).RowHeight = 43.5

# --- Row 22: new "Create shopping cart functionality" entry ----------------
$ws.Range("A22").Value = "ALL"
$ws.Range("B22").Value = "Create shopping cart functionality"
$ws.Range("D22").Value = $false

# --- Row 23: new (in-progress / stub) entry ---------------------------------
$ws.Range("A23").Value = "ALL"
$ws.Range("B23").Value = "Cre"
$ws.Range("D23").Value = $false

# --- Row 11: add a "Solution?" note in column E -----------------------------
$ws.Range("E11").Value = "Also add the option to add payment methods and add fields for this to the database"

# --- Row 12: mark as done ----------------------------------------------------
$ws.Range("D12").Value = $true

# --- Row 19: turn the YouTube URL in column C into a real hyperlink ---------
$ws.Hyperlinks.Add($ws.Range("C19"), "https://www.youtube.com/watch?v=0TnO1GzKWPc&ab_channel=SnippetsCode") | Out-Null
$ws.Rows.Item(19).RowHeight = 43.5

# --- Row 20: mark as done ----------------------------------------------------
$ws.Range("D20").Value = $true

# --- Update the view: selection + scroll position ---------------------------
$win = $excel.Windows.Item(1)
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("E14").Select() | Out-Null
